$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text (remove double space) - AC1, AE1
$ws.Range("AC1").Value = "J6 - Minutos"
$ws.Range("AE1").Value = "J7 - Minutos"

# Update Minutos (K) and J6 - Minutos (AC) values per player row
$ws.Range("K2").Value = 469
$ws.Range("AC2").Value = 45
$ws.Range("K3").Value = 540
$ws.Range("AC3").Value = 90
$ws.Range("K4").Value = 384
$ws.Range("AC4").Value = 25
$ws.Range("K5").Value = 499
$ws.Range("AC5").Value = 90
$ws.Range("K6").Value = 506
$ws.Range("AC6").Value = 85
$ws.Range("K7").Value = 426
$ws.Range("AC7").Value = 90
$ws.Range("K8").Value = 482
$ws.Range("AC8").Value = 90
$ws.Range("AC9").Value = 0
$ws.Range("AC10").Value = 0
$ws.Range("K12").Value = 418
$ws.Range("AC12").Value = 90
$ws.Range("K13").Value = 186
$ws.Range("AC13").Value = 69
$ws.Range("K14").Value = 158
$ws.Range("AC14").Value = 45
$ws.Range("K15").Value = 342
$ws.Range("AC15").Value = 90
$ws.Range("K16").Value = 232
$ws.Range("AC16").Value = 30
$ws.Range("K17").Value = 106
$ws.Range("AC17").Value = 21
$ws.Range("AC18").Value = 0
$ws.Range("AC19").Value = 0
$ws.Range("K21").Value = 150
$ws.Range("AC21").Value = 60
$ws.Range("AC27").Value = 0
$ws.Range("K28").Value = 19
$ws.Range("AC28").Value = 9
